$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

$ws.Range("D3").Value = 93.09999999999999
$ws.Range("C4").Value = 7923
$ws.Range("D4").Value = 97.09999999999999
$ws.Range("C5").Value = 9728
